# "fixed typs in powerpoint"
#
# 1) Slide 7 ("Model Performance"): clarify two bullet points in the
#    accuracy callout textbox.
# 2) Slide 9 ("Future Work"): remove the stray/incomplete "Looking"
#    bullet from the numbered-process SmartArt diagram.

$p = $ppt.ActivePresentation

# --- Slide 7: fix wording in the model-performance bullets -----------------
$s7 = $p.Slides.Item(7)
$statsBox = $s7.Shapes.Item("TextBox 4")
$tr = $statsBox.TextFrame.TextRange

$tr.Replace(
    "Out of  all negative/neutral tweets, we classified 90% of them correctly",
    "Out of  all negative/neutral tweets in our test set, we classified 90% of them correctly",
    0, $false, $false) | Out-Null

$tr.Replace(
    "Of all tweets classified as positive, 70% were actually positive",
    "Of all tweets classified by the model as positive, 70% were actually positive",
    0, $false, $false) | Out-Null

# --- Slide 9: drop the leftover "Looking" SmartArt bullet ------------------
$s9 = $p.Slides.Item(9)
$diagramShape = $s9.Shapes.Item("Content Placeholder 2")
$smartArt = $diagramShape.SmartArt

for ($i = $smartArt.AllNodes.Count; $i -ge 1; $i--) {
    $node = $smartArt.AllNodes.Item($i)
    if ($node.TextFrame2.TextRange.Text -eq "Looking") {
        $node.Delete()
    }
}
